$d = $word.ActiveDocument

# Locate the paragraph that ends with "Hit commit to main"
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd("`r`a`n") -eq "Hit commit to main") {
        $target = $p
        break
    }
}

# Insert a new paragraph right after it, inheriting the list formatting
$rng = $target.Range
$rng.Collapse(0)
$rng.InsertParagraphAfter()

# The newly created paragraph is now the next one; set its text
$newPara = $d.Paragraphs.Item($i + 1)
$newPara.Range.Text = "Find and click the push button"
